# Progress #25 'Interview Preparation': add w:proofErr (spellStart/spellEnd,
# gramStart/gramEnd) markers around flagged words/phrases across the
# document. Word normally inserts these automatically while the user is
# typing / the background proofer runs; here we reproduce the same
# resulting markup by splitting the affected runs and wrapping the
# flagged sub-range with <w:proofErr/> sentinels, using Range.InsertXML
# to replace each whole paragraph with its updated run sequence.

$d = $word.ActiveDocument

function Set-ParagraphXml($AnchorText, $InnerXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($AnchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Set-ParagraphXml: anchor not found: " + $AnchorText)
    }
    $para = $rng.Paragraphs(1).Range

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.InsertXML($xml)
}

# 1) "Is the easiest way to grow your business in 2020" -> gramStart/gramEnd around "2020"
Set-ParagraphXml "Is the easiest way to grow your business in 2020" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">Is the easiest way to grow your business in </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>2020</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
)

# 2) "I want to work at NiceJob" -> spellStart/spellEnd around "NiceJob"
Set-ParagraphXml "I want to work at NiceJob" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">I want to work at </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>NiceJob</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# 3) "For example, when I was working at SiteMax Systems inc. ..." -> spellStart/spellEnd around "SiteMax"
Set-ParagraphXml "For example, when I was working at SiteMax Systems inc." (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">For example, when I was working at </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>SiteMax</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> Systems inc. I was responsible of creating a small front-end application on behalf of Bold construction. I had to dissect the business requirements and deploy solutions after review. It was a long process. But, when it was done, and when a friend of mine who was in charge of delivering web solutions to clients came say &#8220;Hey Moe, Bold construction really loved your solution. So much so that they showed your app to a person at the government.&#8221; I was thrilled by the feeling. I love the feeling of how the solution I created brings smile to customers, and I want to continue this pursuit at this company.</w:t></w:r>'
)

# 4) "...every day of work we do actually helps people..." -> gramStart/gramEnd around "helps"
Set-ParagraphXml "very day of work we do actually helps people, and getting reviews back from business on how " (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>What inspires me is that e</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">very day of work we do actually </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>helps</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> people, and getting reviews back from business on how </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">it has transformed their business is what&#8217;s motivating us to get out from our bed every day </w:t></w:r>'
)

# 5) "One of the requirements of getting a job at companies in south korea ... orindinary Koreans." ->
#    spellStart/spellEnd around "korea" and around "orindinary"
Set-ParagraphXml "One of the requirements of getting a job at companies in south korea" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">One of the requirements of getting a job at companies in south </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>korea</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> is being the best in my field and offer something that&#8217;s not offered by </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>orindinary</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> Koreans.</w:t></w:r>'
)

# 6) "I want to be with my love and I want to work hard to get there" -> gramStart/gramEnd around "love"
Set-ParagraphXml "I want to be with my love and I want to work hard to get there" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">I want to be with my </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>love</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> and I want to work hard to get there</w:t></w:r>'
)

# 7) "I want to stay here at nicejob a long period of time and" (+ remaining runs) -> spellStart/spellEnd around "nicejob"
Set-ParagraphXml "I want to stay here at nicejob a long period of time and" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">I want to stay here at </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>nicejob</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> a long period of time and</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">truly </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>prepare myself to be ready for that moment</w:t></w:r>'
)

# 8) Empty " " list-paragraph (numId 3, between "None" and the bottom-bordered
#    paragraph) loses its numbering/style and becomes a plain indented empty
#    paragraph. Located relative to the unique "None" paragraph since its own
#    text is just a single space.
$rngNone = $d.Content
$rngNone.Find.Execute("None")
$noneParagraph = $rngNone.Paragraphs(1)
$targetParagraph = $noneParagraph.Next()
$targetRange = $targetParagraph.Range
$xmlEmpty = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRange.InsertXML($xmlEmpty)

# 9) "My CTO originally said "Hey moe, don't over do it"." -> spellStart/spellEnd around
#    "moe" and around "over do"
Set-ParagraphXml "originally said" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>My CTO</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> originally said &#8220;Hey </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>moe</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">, don&#8217;t </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>over do</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> it&#8221;.</w:t></w:r>'
)

# 10) "If proportionality constant is  1, then there would be ..." -> gramStart/gramEnd around "is  1"
Set-ParagraphXml "If proportionality constant is  1, then there would be" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">If proportionality constant </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>is  1</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>, then there would be</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> 5/15 zigs and 10/15 zags</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:tab/></w:r>'
)

# 11) "What do you see yourself doing in the next 5 to 10 years" -> gramStart/gramEnd around "years"
Set-ParagraphXml "What do you see yourself doing in the next 5 to 10 years" (
    '<w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:spacing w:val="8"/><w:shd w:val="clear" w:color="auto" w:fill="F9FAFA"/></w:rPr><w:t xml:space="preserve">What do you see yourself doing in the next 5 to 10 </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:spacing w:val="8"/><w:shd w:val="clear" w:color="auto" w:fill="F9FAFA"/></w:rPr><w:t>years</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
)
